# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.705.95"
$ws.Range("E2").Value = "  -0.81%  "

# Row 3
$ws.Range("D3").Value = "2.060.54"
$ws.Range("E3").Value = "  +0.49%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.58%  "

# Row 6
$ws.Range("E6").Value = "  +1.25%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.35"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.46%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.56%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.367"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.01%  "

# Row 11
$ws.Range("E11").Value = "  -2.39%  "

# Row 12
$ws.Range("E12").Value = "  -2.96%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.942"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.79%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.36%  "

# Row 15
$ws.Range("D15").Value = "2.361.48"
$ws.Range("E15").Value = "  +0.42%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.47%  "

# Row 17
$ws.Range("D17").Value = "2.053.02"
$ws.Range("E17").Value = "  -0.83%  "

# Row 18
$ws.Range("D18").Value = "36.630.91"
$ws.Range("E18").Value = "  -0.96%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.19%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.29%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0867"
$ws.Range("E21").Value = "  -2.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.45%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.33%  "

# Row 24
$ws.Range("E24").Value = "  -0.03%  "

# Row 25
$ws.Range("E25").Value = "  -2.35%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.13%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.41%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.05%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.123"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.98%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.50%  "

# Row 32
$ws.Range("E32").Value = "  +7.70%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.54%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0600"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.32%  "

# Row 35
$ws.Range("E35").Value = "  -0.10%  "

# Row 36
$ws.Range("E36").Value = "  -0.26%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0846"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.71%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.12%  "

# Row 39
$ws.Range("E39").Value = "  -4.35%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.05%  "

# Row 41
$ws.Range("E41").Value = "  -4.70%  "

# Row 42
$ws.Range("E42").Value = "  -2.71%  "

# Row 43
$ws.Range("E43").Value = "  -2.97%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "94.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.81%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.15%  "

# Row 46
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.412.00"
$ws.Range("E46").Value = "  +8.68%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0903"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.27%  "

# Row 48
$ws.Range("E48").Value = "  -4.69%  "

# Row 49
$ws.Range("E49").Value = "  +1.36%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.89%  "

# Row 51
$ws.Range("D51").Value = "2.250.35"
$ws.Range("E51").Value = "  +0.62%  "

